$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" share the same data rows; both need the
# "想去人数" (interest count) column F updated for rows 2, 4 and 5.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 247
    $ws.Range("F4").Value = 158
    $ws.Range("F5").Value = 8
}
